# Applies the cryptos-list price/volume refresh described in the commit
# "Updated cryptos list on Mon Oct 30 20:58:11 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "34.486.31"
$ws.Range("E2").Value = "  -0.35%  "

# Row 3
$ws.Range("D3").Value = "1.803.41"
$ws.Range("E3").Value = "  +0.27%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Value = "'228.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.69%  "

# Row 6
$ws.Range("E6").Value = "  +3.96%  "

# Row 7
$ws.Range("E7").Value = "  +0.05%  "

# Row 8
$ws.Range("D8").Value = "'34.86"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.13%  "

# Row 9
$ws.Range("E9").Value = "  +1.28%  "

# Row 11
$ws.Range("E11").Value = "  +0.19%  "

# Row 12
$ws.Range("D12").Value = "2.064.39"
$ws.Range("E12").Value = "  +0.35%  "

# Row 13
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'11.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.90%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.797.00"
$ws.Range("E14").Value = "  +0.25%  "

# Row 15
$ws.Range("D15").Value = "'0.644"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.75%  "

# Row 16
$ws.Range("D16").Value = "34.456.49"

# Row 17
$ws.Range("E17").Value = "  +1.25%  "

# Row 18
$ws.Range("D18").Value = "'69.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.13%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0797"
$ws.Range("E19").Value = "  -0.88%  "

# Row 20
$ws.Range("D20").Value = "'245.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.89%  "

# Row 21
$ws.Range("D21").Value = "'11.52"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.88%  "

# Row 22
$ws.Range("E22").Value = "  +0.13%  "

# Row 23
$ws.Range("E23").Value = "  -0.22%  "

# Row 24
$ws.Range("D24").Value = "'173.82"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.21%  "

# Row 25
$ws.Range("E25").Value = "  +2.09%  "

# Row 26
$ws.Range("D26").Value = "'7.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.04%  "

# Row 27
$ws.Range("E27").Value = "  +1.33%  "

# Row 28
$ws.Range("E28").Value = "  +2.23%  "

# Row 29
$ws.Range("E29").Value = "  -0.14%  "

# Row 30
$ws.Range("E30").Value = "  -1.90%  "

# Row 31
$ws.Range("E31").Value = "  +1.22%  "

# Row 32
$ws.Range("D32").Value = "'3.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.73%  "

# Row 33
$ws.Range("D33").Value = "'1.25"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.29%  "

# Row 34
$ws.Range("D34").Value = "'1.84"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.40%  "

# Row 35
$ws.Range("D35").Value = "'0.685"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.04%  "

# Row 36
$ws.Range("D36").Value = "1.395.09"
$ws.Range("E36").Value = "  -2.34%  "

# Row 37
$ws.Range("E37").Value = "  -2.95%  "

# Row 38
$ws.Range("E38").Value = "  -0.53%  "

# Row 39
$ws.Range("E39").Value = "  -1.08%  "

# Row 40
$ws.Range("D40").Value = "'83.50"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.51%  "

# Row 41
$ws.Range("E41").Value = "  +2.45%  "

# Row 42
$ws.Range("E42").Value = "  +1.45%  "

# Row 44
$ws.Range("E44").Value = "  -0.62%  "

# Row 45
$ws.Range("E45").Value = "  +3.49%  "

# Row 46
$ws.Range("D46").Value = "'0.0510"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.28%  "

# Row 47
$ws.Range("E47").Value = "  -2.20%  "

# Row 48
$ws.Range("D48").Value = "1.963.84"
$ws.Range("E48").Value = "  +0.40%  "

# Row 49
$ws.Range("D49").Value = "'104.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.07%  "

# Row 50
$ws.Range("E50").Value = "  +0.07%  "

# Row 51
$ws.Range("E51").Value = "  +1.30%  "
